$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1050893
$ws.Range("C3").Value = 826821
$ws.Range("C4").Value = 143427
$ws.Range("C5").Value = 65585
$ws.Range("C6").Value = 15060
$ws.Range("C7").Value = 452230
$ws.Range("C8").Value = 598663
$ws.Range("C9").Value = 108050
$ws.Range("C10").Value = 105089
